$wb = $excel.ActiveWorkbook

# With the translation function changed, the upper income limit needs to be
# increased from 200000 to 270000 on the Income sheet (cell B27).
$wsIncome = $wb.Worksheets.Item("Income")
$wsIncome.Range("B27").Value = 270000

# Make the Income sheet the active tab/sheet, with B28 (just below the data)
# selected as the active cell. This also causes the Education sheet (which was
# previously the active tab) to no longer be the selected tab, while keeping
# its own prior selection (B8) intact.
$wsIncome.Activate()
$wsIncome.Range("B28").Select()
